# Stakeholdrers meeting II_Rev2.pptx - apply commit "Update Stakeholdrers meeting II_Rev2.pptx"
#
# Changes:
#  1. Notes master "datum" (date) placeholder text: 22/02/2022 -> 23/02/2022
#  2. Slide 5 ("Project Scope") comparison table:
#     - reposition/resize the table frame
#     - "Hardware available" -> "Hardware partly available"
#     - "Achieve learning goals all members" -> "Achieve learning goals from all members"
#     - "Do not achieve learning goals all members" -> "Do not achieve learning goals from all members"
#     - "Need to understand how to use " -> "Need to understand how to use RT DB"
#     - add new text "Do not achieve the stakeholders' goal" to a previously empty cell

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Notes master date placeholder
# ---------------------------------------------------------------------------
$nm = $p.NotesMaster
$dateFooter = $nm.HeadersFooters.DateAndTime
$dateFooter.Text = "23/02/2022"

# ---------------------------------------------------------------------------
# 2. Slide 5 table
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)

# Reposition / resize the graphic frame (EMU 1550404,855429 6385243x3769360).
# Shape.Left/Top/Width/Height are expressed in points (1 pt = 12700 EMU); the
# literal values below are chosen so the round-trip lands on the exact EMU
# target despite the single-precision float storage used internally.
$tableShape.Left = 122.07905578613281
$tableShape.Top = 67.35662078857422
$tableShape.Height = 296.8000183105469
# Width (cx) is unchanged from the original 6385243 EMU, so it is left alone.

$tbl = $tableShape.Table

# "Hardware available" -> "Hardware partly available"
$tbl.Cell(2, 4).Shape.TextFrame.TextRange.Text = "Hardware partly available"

# "Achieve learning goals all members" -> "Achieve learning goals from all members"
$tbl.Cell(3, 4).Shape.TextFrame.TextRange.Text = "Achieve learning goals from all members"

# "Do not achieve learning goals all members" -> "Do not achieve learning goals from all members"
$tbl.Cell(5, 2).Shape.TextFrame.TextRange.Text = "Do not achieve learning goals from all members"

# "Need to understand how to use " -> "Need to understand how to use RT DB"
$tbl.Cell(6, 3).Shape.TextFrame.TextRange.Text = "Need to understand how to use RT DB"

# Previously-empty cell gets a new sentence (curly apostrophe, U+2019).
$apos = [char]0x2019
$tbl.Cell(7, 2).Shape.TextFrame.TextRange.Text = "Do not achieve the stakeholders" + $apos + " goal"
